$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values for the rows whose price changed
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.67"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.98"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.401"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05937"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.396"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8074"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9131"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1419"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07438"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03336"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03074"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09344"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.937"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001574"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04797"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005936"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005476"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004447"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009860"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00007501"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.665"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3247"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1347"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002446"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03898"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006205"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002743"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006596"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005199"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005796"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.048"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"

# Update Hora (column G) values for all data rows (14 -> 15)
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "15"
